$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.411.13"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.629.32"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.60"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.80%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3641"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -2.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08209"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.229"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.36"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.531"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001247"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.14%  "

$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.624.71"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.95"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06970"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.67"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.527"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.65"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.400.41"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.141"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.448"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.17%  "

$ws.Range("E27").Value = "  +0.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.58"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.297"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.14%  "

$ws.Range("E30").Value = "  -2.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.800.34"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.243"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.801"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.039"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.81"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02781"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2513"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08742"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07100"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.985"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7034"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.346"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.09"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6546"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.287"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.971"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08003"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.46%  "

$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.86"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.36%  "
